$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.069.10"
$ws.Range("E2").Value = "  -0.17%  "
$ws.Range("D3").Value = "1.621.85"
$ws.Range("E3").Value = "  -0.89%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.11"
$ws.Range("E5").Value = "  -1.23%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.0627"
$ws.Range("E8").Value = "  +0.45%  "
$ws.Range("E9").Value = "  -1.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.90"
$ws.Range("E10").Value = "  -0.09%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0840"
$ws.Range("E11").Value = "  -0.62%  "
$ws.Range("D12").Value = "1.848.96"
$ws.Range("E12").Value = "  -0.88%  "
$ws.Range("D13").Value = "1.626.79"
$ws.Range("E13").Value = "  -0.59%  "
$ws.Range("E14").Value = "  -0.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.538"
$ws.Range("E15").Value = "  -0.36%  "
$ws.Range("D16").Value = "27.048.67"
$ws.Range("E16").Value = "  -0.25%  "
$ws.Range("E17").Value = "  -3.26%  "
$ws.Range("D18").Value = "0.0₃0736"
$ws.Range("E18").Value = "  -0.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "214.36"
$ws.Range("E19").Value = "  -1.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.82"
$ws.Range("E21").Value = "  -0.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.33"
$ws.Range("E22").Value = "  -1.74%  "
$ws.Range("E23").Value = "  -7.82%  "
$ws.Range("E24").Value = "  -0.97%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.52"
$ws.Range("E25").Value = "  +0.63%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.43"
$ws.Range("E26").Value = "  +0.55%  "
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("E28").Value = "  -3.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.49"
$ws.Range("E29").Value = "  -1.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0510"
$ws.Range("E30").Value = "  +0.76%  "
$ws.Range("E31").Value = "  -1.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.31"
$ws.Range("E32").Value = "  -1.87%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.730"
$ws.Range("E33").Value = "  +34.70%  "
$ws.Range("E34").Value = "  -0.46%  "
$ws.Range("D35").Value = "1.335.06"
$ws.Range("E35").Value = "  +2.57%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.55"
$ws.Range("E36").Value = "  -1.06%  "
$ws.Range("E37").Value = "  -0.69%  "
$ws.Range("E38").Value = "  -0.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.840"
$ws.Range("E39").Value = "  -1.67%  "
$ws.Range("E40").Value = "  -0.13%  "
$ws.Range("E41").Value = "  -0.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.795"
$ws.Range("E42").Value = "  -1.46%  "
$ws.Range("E43").Value = "  +0.61%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.94"
$ws.Range("E44").Value = "  +3.61%  "
$ws.Range("D45").Value = "1.760.36"
$ws.Range("E45").Value = "  -0.89%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "89.83"
$ws.Range("E46").Value = "  -1.51%  "
$ws.Range("E47").Value = "  +2.67%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.853"
$ws.Range("E48").Value = "  +27.92%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0999"
$ws.Range("E49").Value = "  +4.28%  "
$ws.Range("E50").Value = "  -0.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.59"
$ws.Range("E51").Value = "  -0.92%  "
